$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "29.300.37"
$ws.Range("E2").Value = "  +0.42%  "

Set-TextValue $ws "D3" "1.875.07"
$ws.Range("E3").Value = "  +0.69%  "

Set-TextValue $ws "D4" "0.9997"
$ws.Range("E4").Value = "  -0.14%  "

Set-TextValue $ws "D5" "0.7135"
$ws.Range("E5").Value = "  -0.76%  "

Set-TextValue $ws "D6" "241.92"
$ws.Range("E6").Value = "  +0.56%  "

Set-TextValue $ws "D7" "0.9999"
$ws.Range("E7").Value = "  -0.19%  "

Set-TextValue $ws "D8" "0.3110"
$ws.Range("E8").Value = "  +1.26%  "

Set-TextValue $ws "D9" "0.07729"
$ws.Range("E9").Value = "  -0.01%  "

Set-TextValue $ws "D10" "25.11"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("E11").Value = "  +1.57%  "

Set-TextValue $ws "D12" "1.880.79"
$ws.Range("E12").Value = "  -0.26%  "

Set-TextValue $ws "D13" "5.213"
$ws.Range("E13").Value = "  +0.18%  "

Set-TextValue $ws "D14" "0.7113"
$ws.Range("E14").Value = "  -0.46%  "

Set-TextValue $ws "D15" "91.30"
$ws.Range("E15").Value = "  +1.38%  "

Set-TextValue $ws "D16" "29.286.71"
$ws.Range("E16").Value = "  +0.27%  "

Set-TextValue $ws "D17" "0.000008297"
$ws.Range("E17").Value = "  +6.73%  "

Set-TextValue $ws "D18" "5.988"
$ws.Range("E18").Value = "  +2.97%  "

Set-TextValue $ws "D19" "242.72"
$ws.Range("E19").Value = "  +0.04%  "

Set-TextValue $ws "D20" "2.127.97"
$ws.Range("E20").Value = "  +0.41%  "

Set-TextValue $ws "D21" "13.21"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("E22").Value = "  -0.23%  "

Set-TextValue $ws "D23" "7.808"
$ws.Range("E23").Value = "  -1.65%  "

Set-TextValue $ws "D24" "0.9998"
$ws.Range("E24").Value = "  -0.22%  "

Set-TextValue $ws "D25" "0.1620"
$ws.Range("E25").Value = "  +2.03%  "

Set-TextValue $ws "D26" "163.13"
$ws.Range("E26").Value = "  +0.69%  "

Set-TextValue $ws "D27" "9.018"
$ws.Range("E27").Value = "  +1.43%  "

$ws.Range("E28").Value = "  +2.07%  "

Set-TextValue $ws "D29" "1.504"
$ws.Range("E29").Value = "  +0.84%  "

Set-TextValue $ws "D30" "4.414"
$ws.Range("E30").Value = "  +1.54%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D31" "4.324"
$ws.Range("E31").Value = "  +5.97%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D32" "1.286"
$ws.Range("E32").Value = "  -0.80%  "

Set-TextValue $ws "D33" "0.05257"
$ws.Range("E33").Value = "  +1.58%  "

Set-TextValue $ws "D34" "1.924"
$ws.Range("E34").Value = "  +0.73%  "

$ws.Range("E35").Value = "  -0.01%  "

Set-TextValue $ws "D36" "0.7469"
$ws.Range("E36").Value = "  +2.72%  "

Set-TextValue $ws "D37" "2.684"
$ws.Range("E37").Value = "  +0.16%  "

Set-TextValue $ws "D38" "0.01859"
$ws.Range("E38").Value = "  +0.74%  "

Set-TextValue $ws "D39" "2.720"
$ws.Range("E39").Value = "  +1.03%  "

Set-TextValue $ws "D40" "1.155.31"
$ws.Range("E40").Value = "  -0.12%  "

Set-TextValue $ws "D41" "6.365"
$ws.Range("E41").Value = "  +4.71%  "

Set-TextValue $ws "D42" "73.12"
$ws.Range("E42").Value = "  +1.58%  "

Set-TextValue $ws "D43" "0.8849"
$ws.Range("E43").Value = "  -1.64%  "

Set-TextValue $ws "D44" "106.21"
$ws.Range("E44").Value = "  +4.59%  "

Set-TextValue $ws "D45" "0.9995"
$ws.Range("E45").Value = "  -0.24%  "

Set-TextValue $ws "D46" "2.024.60"
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("E47").Value = "  +2.62%  "

Set-TextValue $ws "D48" "0.5192"
$ws.Range("E48").Value = "  -1.66%  "

Set-TextValue $ws "D49" "9.388"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("E50").Value = "  +2.92%  "

Set-TextValue $ws "D51" "0.4300"
$ws.Range("E51").Value = "  +1.77%  "
